$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header "time_taken" in F1 and give it the same style as the
# other header cells (e.g. E1 - bold, bordered, centered style).
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)

# Populate F2:F15 with the recorded time_taken timestamps (as plain text,
# matching the inlineStr/text cells from the source data).
$timestamps = @(
    "2021-10-05 13:39:04.009859",
    "2021-10-05 13:39:04.009873",
    "2021-10-05 13:39:04.009877",
    "2021-10-05 13:39:04.009881",
    "2021-10-05 13:39:04.009884",
    "2021-10-05 13:39:04.009887",
    "2021-10-05 13:39:04.009890",
    "2021-10-05 13:39:04.009894",
    "2021-10-05 13:39:04.009897",
    "2021-10-05 13:39:04.009900",
    "2021-10-05 13:39:04.009903",
    "2021-10-05 13:39:04.009906",
    "2021-10-05 13:39:04.009909",
    "2021-10-05 13:39:04.009912"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
